$d = $word.ActiveDocument

# --- Change 1 ---
# Four adjacent runs with identical formatting are combined into a single
# run. The wording itself does not change, only the run it lives in, so we
# re-apply the exact same text across the full span to coalesce the runs.
$c1Old = " This is a red flag that Louise did not have many donors. Louise had only 10 donors. The average number of donors of successful play kickstarters in her range was 44. That is less than a fourth of the average number of donors for successful plays. In fact, it looks like average number of donors for failed campaigns across all buckets is 8, while successful plays have an average of 56 donors. Hence, the fewer number of donors there are for the plays in any of the bucket range, the greater the chances they will fail."
$d.Content.Find.Execute($c1Old, $true, $false, $false, $false, $false, $true, 1, $false, $c1Old, 2) | Out-Null

# --- Change 2 ---
# "The formula and table was fairly simple" -> "...table were fairly simple"
$d.Content.Find.Execute("table was fairly", $true, $false, $false, $false, $false, $true, 1, $false, "table were fairly", 2) | Out-Null

# --- Change 3 ---
# "Louise's play fever since" -> "Louise's play, Fever, since"
$d.Content.Find.Execute("play fever since", $true, $false, $false, $false, $false, $true, 1, $false, "play, Fever, since", 2) | Out-Null
